# DTProgram.xlsx — "fixed validate input & modifie read config input key field"
#
# The lookup/config table on the active sheet ("Sheet1 (2)") used the row
# IDs "C2"/"C3" for the two conditional-threshold rows; the key field for
# those rows is renamed to "X2"/"X3" (new entries appended to the shared
# string table, matching how Excel itself grows sst rather than mutating
# the old "C2"/"C3" strings, which are still referenced elsewhere in the
# sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 and row 4 key-field (column A) IDs: C2 -> X2, C3 -> X3
$ws.Range("A3").Value = "X2"
$ws.Range("A4").Value = "X3"

# Column A widened to comfortably fit the (now slightly different) key
# values.
$ws.Columns.Item(1).ColumnWidth = 6.5
